# Insert a new weekly price record as row 60 on the single worksheet,
# shifting all subsequent rows down by one (dimension grows from
# A1:R150 to A1:R151).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 60 (and everything below it) down by one row.
$ws.Rows(60).Insert()

# Populate the newly inserted row 60 with the new record.
$ws.Cells.Item(60, 1).Value2  = 4
$ws.Cells.Item(60, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(60, 3).Value2  = "Los Lagos"
$ws.Cells.Item(60, 4).Value2  = 44477
$ws.Cells.Item(60, 5).Value2  = 10
$ws.Cells.Item(60, 6).Value2  = 100112021
$ws.Cells.Item(60, 7).Value2  = "Ají"
$ws.Cells.Item(60, 8).Value2  = "Inferno"
$ws.Cells.Item(60, 9).Value2  = "Primera"
$ws.Cells.Item(60, 10).Value2 = 140
$ws.Cells.Item(60, 11).Value2 = 48000
$ws.Cells.Item(60, 12).Value2 = 50000
$ws.Cells.Item(60, 13).Value2 = 49000
$ws.Cells.Item(60, 14).Value2 = "$/caja 12 kilos"
$ws.Cells.Item(60, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(60, 16).Value2 = 4083
$ws.Cells.Item(60, 17).Value2 = 12
$ws.Cells.Item(60, 18).Value2 = "Hortaliza"
